$wb = $excel.ActiveWorkbook

# Sheet 1: "Metadata"
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$ws1.Range("B5").Value = "CDACompressionAlgorithm"
$ws1.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$ws1.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"
